# leitor_cnpj/cnpj_ler.xlsx — fill in the remaining CNPJ rows (A3:A5) and
# append a new CNPJ on row 6, matching the formatting already used for A2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 1612155000141
$ws.Range("A4").Value = 28129260000938
$ws.Range("A5").Value = 52123916001457

# New row 6 needs to carry the same formatting (number format / font) as the
# row above it, so copy A5's formatting down into A6 before writing the value.
$ws.Range("A5").Copy($ws.Range("A6"))
$ws.Rows.Item(6).RowHeight = $ws.Rows.Item(5).RowHeight

$ws.Range("A6").Value = 89425888000541

$ws.Range("A6").Select()
